# Update "想去人数" (column F) counts across the 展览, 演出 and 全部类型 sheets
# to reflect the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value  = 123
$ws.Cells.Item(3, 6).Value  = 729
$ws.Cells.Item(4, 6).Value  = 53
$ws.Cells.Item(6, 6).Value  = 2990
$ws.Cells.Item(8, 6).Value  = 2007
$ws.Cells.Item(10, 6).Value = 299
$ws.Cells.Item(11, 6).Value = 863
$ws.Cells.Item(12, 6).Value = 951
$ws.Cells.Item(13, 6).Value = 200
$ws.Cells.Item(14, 6).Value = 420
$ws.Cells.Item(17, 6).Value = 65
$ws.Cells.Item(19, 6).Value = 7235
$ws.Cells.Item(20, 6).Value = 267
$ws.Cells.Item(21, 6).Value = 1960
$ws.Cells.Item(22, 6).Value = 188
$ws.Cells.Item(24, 6).Value = 163
$ws.Cells.Item(25, 6).Value = 447
$ws.Cells.Item(26, 6).Value = 500
$ws.Cells.Item(27, 6).Value = 75
$ws.Cells.Item(29, 6).Value = 947
$ws.Cells.Item(31, 6).Value = 121
$ws.Cells.Item(33, 6).Value = 1121
$ws.Cells.Item(34, 6).Value = 1915
$ws.Cells.Item(35, 6).Value = 475
$ws.Cells.Item(36, 6).Value = 12
$ws.Cells.Item(38, 6).Value = 254
$ws.Cells.Item(40, 6).Value = 150
$ws.Cells.Item(41, 6).Value = 276
$ws.Cells.Item(43, 6).Value = 195

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 9
$ws.Cells.Item(7, 6).Value = 5

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 123
$ws.Cells.Item(3, 6).Value  = 729
$ws.Cells.Item(4, 6).Value  = 53
$ws.Cells.Item(9, 6).Value  = 2990
$ws.Cells.Item(11, 6).Value = 2007
$ws.Cells.Item(13, 6).Value = 299
$ws.Cells.Item(14, 6).Value = 863
$ws.Cells.Item(16, 6).Value = 951
$ws.Cells.Item(17, 6).Value = 200
$ws.Cells.Item(18, 6).Value = 420
$ws.Cells.Item(21, 6).Value = 65
$ws.Cells.Item(23, 6).Value = 7236
$ws.Cells.Item(24, 6).Value = 267
$ws.Cells.Item(25, 6).Value = 1961
$ws.Cells.Item(26, 6).Value = 9
$ws.Cells.Item(27, 6).Value = 188
$ws.Cells.Item(29, 6).Value = 163
$ws.Cells.Item(30, 6).Value = 447
$ws.Cells.Item(31, 6).Value = 500
$ws.Cells.Item(32, 6).Value = 75
$ws.Cells.Item(34, 6).Value = 947
$ws.Cells.Item(36, 6).Value = 121
$ws.Cells.Item(37, 6).Value = 1121
$ws.Cells.Item(38, 6).Value = 1915
$ws.Cells.Item(39, 6).Value = 475
$ws.Cells.Item(40, 6).Value = 12
$ws.Cells.Item(42, 6).Value = 254
$ws.Cells.Item(44, 6).Value = 150
$ws.Cells.Item(45, 6).Value = 276
$ws.Cells.Item(46, 6).Value = 5
$ws.Cells.Item(49, 6).Value = 195
